$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- xlPasteFormats constant ---
$xlPasteFormats = -4122

# --- Update existing cell text (rows whose A/B/C layout is unchanged) ---
$ws.Range("B13").Value = '5840841 - Gilberto Garcia Cortez'
$ws.Range("C13").Value = '5840841 - Gilberto Garcia Cortez'

$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = '1- Introdução:2- Coeficiente de difusão:3- Concentrações, velocidade e fluxos:4 -Equações da continuidade em transferência de massa:5- Difusão em regime permanente sem reação química:6- Difusão com reação química:7- Transferência de massa entre fases.'
$ws.Range("C14").Value = '1- Introdução:2- Coeficiente de difusão:3- Concentrações, velocidade e fluxos:4 -Equações da continuidade em transferência de massa:5- Difusão em regime permanente sem reação química:6- Difusão com reação química:7- Transferência de massa entre fases.'

$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = '1 - Introduction: 2 - Diffusion coefficient: 3 - Concentrations, and flow rate: 4 - Equation of continuity for mass transfer: 5 - Diffusion in continuous operation without chemical reaction: 6 - Diffusion with chemical reaction: 7 - Mass transfer between phases.'
$ws.Range("C15").Value = '1 - Introduction: 2 - Diffusion coefficient: 3 - Concentrations, and flow rate: 4 - Equation of continuity for mass transfer: 5 - Diffusion in continuous operation without chemical reaction: 6 - Diffusion with chemical reaction: 7 - Mass transfer between phases.'

$ws.Range("A17").Value = 'Syllabus:'

$ws.Range("A18").Value = 'Avaliação:'

$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Range("C19").Value = 'A avaliação será feita por meio de provas escritas.'

$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'A Nota Final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3.'
$ws.Range("C20").Value = 'A Nota Final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3.'

$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma Prova Escrita (PE) e a Média de Recuperação (MR) será calculada pela fórmula: MR = (NF + PE)/2.'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma Prova Escrita (PE) e a Média de Recuperação (MR) será calculada pela fórmula: MR = (NF + PE)/2.'

$ws.Range("A22").Value = 'Bibliografia:'

# --- Remove cells that no longer exist in the new layout ---
$ws.Range("A13").Clear()
$ws.Range("B18:C18").Clear()
$ws.Range("B23:C23").Clear()

# --- Add new cells, copying number/alignment/font formats from a donor cell ---
$ws.Range("B16:C16").Copy()
$ws.Range("B17:C17").PasteSpecial($xlPasteFormats)
$ws.Range("B17").Value = '1 - Introduction: Mass transfer: Definition. Classification of operations involving mass transfer. Contributions to mass transfer. Types of diffusion. 2 - Diffusion coefficient: Considerations. Diffusion in gases: Analysis of Fick''s first law, the diffusion coefficient for gases. Estimation of the diffusion coefficient from a known diffusion coefficient at high temperature and pressure. Diffusion coefficient of a solute in a multicomponent mixture of stagnant gases. Diffusion in liquids. Diffusion in solids. 3 - Concentrations, and flow rate: Concentration. Speed and Flow. The equation of Stefan – Maxwell. 4 - Equation of continuity for mass transfer: Equations of continuity of a solute molar. Transient without/with speed zero means. Medium with and without chemical reaction. 5 - Diffusion in permanent regime without chemical reaction: One-dimensional diffusion in permanent regime. Diffusion through inert and stagnant gaseous film. Pseudo-stationary diffusion in a stagnant gaseous film. Equimolar contradifusion. Molar rate in isolated beads. Diffusion in membranes. 6 - Diffusion in permanent with chemical reaction: Diffusion in steady state with heterogeneous chemical reaction on the surface of a nonporous catalytic particle. Diffusion with heterogeneous chemical reaction on the surface of a non-catalytic and non-porous particle. Intraparticle diffusion with heterogeneous chemical reaction. Continuous diffusion with homogeneous chemical reaction. 7 - Mass transfer between phases: Theory of the two resistors. Individual and global coefficient of mass transfer. Global mass transfer coefficients. Volumetric coefficients of transfer of mass to towers of fillings. Balance macroscope of matter. Continuous operations.'
$ws.Range("C17").Value = '1 - Introduction: Mass transfer: Definition. Classification of operations involving mass transfer. Contributions to mass transfer. Types of diffusion. 2 - Diffusion coefficient: Considerations. Diffusion in gases: Analysis of Fick''s first law, the diffusion coefficient for gases. Estimation of the diffusion coefficient from a known diffusion coefficient at high temperature and pressure. Diffusion coefficient of a solute in a multicomponent mixture of stagnant gases. Diffusion in liquids. Diffusion in solids. 3 - Concentrations, and flow rate: Concentration. Speed and Flow. The equation of Stefan – Maxwell. 4 - Equation of continuity for mass transfer: Equations of continuity of a solute molar. Transient without/with speed zero means. Medium with and without chemical reaction. 5 - Diffusion in permanent regime without chemical reaction: One-dimensional diffusion in permanent regime. Diffusion through inert and stagnant gaseous film. Pseudo-stationary diffusion in a stagnant gaseous film. Equimolar contradifusion. Molar rate in isolated beads. Diffusion in membranes. 6 - Diffusion in permanent with chemical reaction: Diffusion in steady state with heterogeneous chemical reaction on the surface of a nonporous catalytic particle. Diffusion with heterogeneous chemical reaction on the surface of a non-catalytic and non-porous particle. Intraparticle diffusion with heterogeneous chemical reaction. Continuous diffusion with homogeneous chemical reaction. 7 - Mass transfer between phases: Theory of the two resistors. Individual and global coefficient of mass transfer. Global mass transfer coefficients. Volumetric coefficients of transfer of mass to towers of fillings. Balance macroscope of matter. Continuous operations.'

$ws.Range("B16:C16").Copy()
$ws.Range("B22:C22").PasteSpecial($xlPasteFormats)
$ws.Range("B22").Value = '1) CREMASCO, M. A. Fundamentos de Transferência de Massa. Campinas: Editora Unicamp; 2008-2009.2) INCROPERA, F. P.; WITT, D. P. Fundamentos de Transferência de Calor e Massa. Rio de Janeiro: LTC, 2008.3) Bird, R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. 2ª ed. Rio de Janeiro: LTC, 2004.4) BENNETT, C. O.; MYERS, J. E. Fenômeno de Transporte: Quantidade de Movimento, Calor e Massa. São Paulo: McGrawc- Hill, 1978.5) COULSON, J. M.; RICHARDSON, J. F.; BACKHURST, J. R.; HARKER, J. H. Fluid Flow, Heat Transfer and Mass Transfer. In: COULSON & Richardson Series - Chemical Engineering. 5th ed. Pergamon Press, Oxford, 1996. v.16) FOUST, A. S.; Wenzel, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSON, L. B. Princípios das Operações Unitárias. 2ª ed. Rio de Janeiro: Guanabara Dois, 1982.7) PERRY''s Chemical Engineers Handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry. New York: McGraw-Hill, 2008.8) WELTY, J. R.; PIGFORD, R. L.; WILKE, C. R. Fundamentals of Momentum, Heat, and Mass Transfer. 5th ed. USA: John Wiley & Sons, Inc, 2008.9) POLING, B. E.; PRAUSNITZ, J. M.; O''CONNELL, J. The Properties of Gases and Liquids. 5th ed. New York: McGraw-Hill, 2000.10) CALDAS, J. N.; DE LACERDA, A. I.; VELOSO, E.; PASCHOAL, L. C. M. Internos de Torres: Pratos & Recheios. 2ª ed. Rio de Janeiro: Editora Interciência, 2007.'
$ws.Range("C22").Value = '1) CREMASCO, M. A. Fundamentos de Transferência de Massa. Campinas: Editora Unicamp; 2008-2009.2) INCROPERA, F. P.; WITT, D. P. Fundamentos de Transferência de Calor e Massa. Rio de Janeiro: LTC, 2008.3) Bird, R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. 2ª ed. Rio de Janeiro: LTC, 2004.4) BENNETT, C. O.; MYERS, J. E. Fenômeno de Transporte: Quantidade de Movimento, Calor e Massa. São Paulo: McGrawc- Hill, 1978.5) COULSON, J. M.; RICHARDSON, J. F.; BACKHURST, J. R.; HARKER, J. H. Fluid Flow, Heat Transfer and Mass Transfer. In: COULSON & Richardson Series - Chemical Engineering. 5th ed. Pergamon Press, Oxford, 1996. v.16) FOUST, A. S.; Wenzel, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSON, L. B. Princípios das Operações Unitárias. 2ª ed. Rio de Janeiro: Guanabara Dois, 1982.7) PERRY''s Chemical Engineers Handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry. New York: McGraw-Hill, 2008.8) WELTY, J. R.; PIGFORD, R. L.; WILKE, C. R. Fundamentals of Momentum, Heat, and Mass Transfer. 5th ed. USA: John Wiley & Sons, Inc, 2008.9) POLING, B. E.; PRAUSNITZ, J. M.; O''CONNELL, J. The Properties of Gases and Liquids. 5th ed. New York: McGraw-Hill, 2000.10) CALDAS, J. N.; DE LACERDA, A. I.; VELOSO, E.; PASCHOAL, L. C. M. Internos de Torres: Pratos & Recheios. 2ª ed. Rio de Janeiro: Editora Interciência, 2007.'

$ws.Range("A18").Copy()
$ws.Range("A23").PasteSpecial($xlPasteFormats)
$ws.Range("A23").Value = 'Requisitos:'

$ws.Range("B10:C10").Copy()
$ws.Range("B24:C24").PasteSpecial($xlPasteFormats)
$ws.Range("B24").Value = 'LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)
'

$excel.CutCopyMode = $false

# --- Row heights ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(24).RowHeight = 30

